$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 4
$ws.Range("K4").Value = 9
$ws.Range("K5").Value = 9
$ws.Range("K7").Value = 9
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 9
$ws.Range("G12").Value = 13
$ws.Range("H12").Value = 10
$ws.Range("J12").Value = 20
$ws.Range("E15").Value = 19
$ws.Range("D16").Value = 35
$ws.Range("E16").Value = 19
$ws.Range("F16").Value = 24
$ws.Range("K18").Value = 9
$ws.Range("F23").Value = 24
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 10
$ws.Range("I24").Value = 20
$ws.Range("K24").Value = 9
$ws.Range("D25").Value = 40
$ws.Range("K25").Value = 9
$ws.Range("E38").Value = 19
$ws.Range("G38").Value = 13
$ws.Range("H38").Value = 10

$ws.Range("K10").Select()
